$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.484391333333334
$ws.Cells.Item(2, 8).Value = 4.453174000000001
$ws.Cells.Item(2, 9).Value = 0.14812258302985
$ws.Cells.Item(2, 10).Value = 0.157784627403132
$ws.Cells.Item(2, 13).Value = 36.48539666666667
$ws.Cells.Item(2, 14).Value = 109.45619
$ws.Cells.Item(2, 15).Value = 0.4260639713374229
$ws.Cells.Item(2, 16).Value = 0.4324607845540777
$ws.Cells.Item(2, 17).Value = 54.1586066052289
$ws.Cells.Item(2, 18).Value = 487.4274594470601
$ws.Cells.Item(2, 19).Value = 0.06310969597045504
$ws.Cells.Item(2, 20).Value = 0.06823566375733128
$ws.Cells.Item(3, 7).Value = 1.484391333333334
$ws.Cells.Item(3, 8).Value = 4.453174000000001
$ws.Cells.Item(3, 9).Value = 0.14812258302985
$ws.Cells.Item(3, 10).Value = 0.157784627403132
$ws.Cells.Item(3, 15).Value = 0.1743777127077069
$ws.Cells.Item(3, 16).Value = 0.1769957741547643
$ws.Cells.Item(3, 17).Value = 22.16581212819112
$ws.Cells.Item(3, 18).Value = 199.49230915372
$ws.Cells.Item(3, 19).Value = 0.02582927722910264
$ws.Cells.Item(3, 20).Value = 0.02792721227693839
$ws.Cells.Item(4, 7).Value = 1.484391333333334
$ws.Cells.Item(4, 8).Value = 4.453174000000001
$ws.Cells.Item(4, 9).Value = 0.14812258302985
$ws.Cells.Item(4, 10).Value = 0.157784627403132
$ws.Cells.Item(4, 13).Value = 9.680823666666667
$ws.Cells.Item(4, 14).Value = 29.042471
$ws.Cells.Item(4, 15).Value = 0.1130493445068016
$ws.Cells.Item(4, 16).Value = 0.1147466378470605
$ws.Cells.Item(4, 17).Value = 14.37013075032823
$ws.Cells.Item(4, 18).Value = 129.331176752954
$ws.Cells.Item(4, 19).Value = 0.01674516091817884
$ws.Cells.Item(4, 20).Value = 0.01810525549846057
$ws.Cells.Item(5, 7).Value = 1.484391333333334
$ws.Cells.Item(5, 8).Value = 4.453174000000001
$ws.Cells.Item(5, 9).Value = 0.14812258302985
$ws.Cells.Item(5, 10).Value = 0.157784627403132
$ws.Cells.Item(5, 13).Value = 3.79999
$ws.Cells.Item(5, 14).Value = 7.59998
$ws.Cells.Item(5, 15).Value = 0.04437498227672168
$ws.Cells.Item(5, 16).Value = 0.0300274777826206
$ws.Cells.Item(5, 17).Value = 5.640672222753334
$ws.Cells.Item(5, 18).Value = 33.84403333652001
$ws.Cells.Item(5, 19).Value = 0.006572936996731829
$ws.Cells.Item(5, 20).Value = 0.004737874393786615
$ws.Cells.Item(6, 7).Value = 1.484391333333334
$ws.Cells.Item(6, 8).Value = 4.453174000000001
$ws.Cells.Item(6, 9).Value = 0.14812258302985
$ws.Cells.Item(6, 10).Value = 0.157784627403132
$ws.Cells.Item(6, 13).Value = 20.734808
$ws.Cells.Item(6, 14).Value = 62.204424
$ws.Cells.Item(6, 15).Value = 0.242133989171347
$ws.Cells.Item(6, 16).Value = 0.245769325661477
$ws.Cells.Item(6, 17).Value = 30.77856929353067
$ws.Cells.Item(6, 18).Value = 277.007123641776
$ws.Cells.Item(6, 19).Value = 0.03586551191538165
$ws.Cells.Item(6, 20).Value = 0.03877862147661515
$ws.Cells.Item(7, 9).Value = 0.4232592596904365
$ws.Cells.Item(7, 10).Value = 0.4508684848665014
$ws.Cells.Item(7, 13).Value = 36.48539666666667
$ws.Cells.Item(7, 14).Value = 109.45619
$ws.Cells.Item(7, 15).Value = 0.4260639713374229
$ws.Cells.Item(7, 16).Value = 0.4324607845540777
$ws.Cells.Item(7, 17).Value = 154.7578449464067
$ws.Cells.Item(7, 18).Value = 1392.82060451766
$ws.Cells.Item(7, 19).Value = 0.180335521089045
$ws.Cells.Item(7, 20).Value = 0.1949829386960755
$ws.Cells.Item(8, 9).Value = 0.4232592596904365
$ws.Cells.Item(8, 10).Value = 0.4508684848665014
$ws.Cells.Item(8, 15).Value = 0.1743777127077069
$ws.Cells.Item(8, 16).Value = 0.1769957741547643
$ws.Cells.Item(8, 19).Value = 0.07380698158717564
$ws.Cells.Item(8, 20).Value = 0.07980181652093206
$ws.Cells.Item(9, 9).Value = 0.4232592596904365
$ws.Cells.Item(9, 10).Value = 0.4508684848665014
$ws.Cells.Item(9, 13).Value = 9.680823666666667
$ws.Cells.Item(9, 14).Value = 29.042471
$ws.Cells.Item(9, 15).Value = 0.1130493445068016
$ws.Cells.Item(9, 16).Value = 0.1147466378470605
$ws.Cells.Item(9, 17).Value = 41.06254953583267
$ws.Cells.Item(9, 18).Value = 369.562945822494
$ws.Cells.Item(9, 19).Value = 0.04784918186443798
$ws.Cells.Item(9, 20).Value = 0.05173564274962933
$ws.Cells.Item(10, 9).Value = 0.4232592596904365
$ws.Cells.Item(10, 10).Value = 0.4508684848665014
$ws.Cells.Item(10, 13).Value = 3.79999
$ws.Cells.Item(10, 14).Value = 7.59998
$ws.Cells.Item(10, 15).Value = 0.04437498227672168
$ws.Cells.Item(10, 16).Value = 0.0300274777826206
$ws.Cells.Item(10, 17).Value = 16.11818198362
$ws.Cells.Item(10, 18).Value = 96.70909190172
$ws.Cells.Item(10, 19).Value = 0.01878212214722146
$ws.Cells.Item(10, 20).Value = 0.01353844341221268
$ws.Cells.Item(11, 9).Value = 0.4232592596904365
$ws.Cells.Item(11, 10).Value = 0.4508684848665014
$ws.Cells.Item(11, 13).Value = 20.734808
$ws.Cells.Item(11, 14).Value = 62.204424
$ws.Cells.Item(11, 15).Value = 0.242133989171347
$ws.Cells.Item(11, 16).Value = 0.245769325661477
$ws.Cells.Item(11, 17).Value = 87.949549535504
$ws.Cells.Item(11, 18).Value = 791.5459458195361
$ws.Cells.Item(11, 19).Value = 0.1024854530025565
$ws.Cells.Item(11, 20).Value = 0.1108096434876519
$ws.Cells.Item(12, 7).Value = 0.8171586666666667
$ws.Cells.Item(12, 8).Value = 2.451476
$ws.Cells.Item(12, 9).Value = 0.08154160546066344
$ws.Cells.Item(12, 10).Value = 0.08686056894424524
$ws.Cells.Item(12, 13).Value = 36.48539666666667
$ws.Cells.Item(12, 14).Value = 109.45619
$ws.Cells.Item(12, 15).Value = 0.4260639713374229
$ws.Cells.Item(12, 16).Value = 0.4324607845540777
$ws.Cells.Item(12, 17).Value = 29.81435809293778
$ws.Cells.Item(12, 18).Value = 268.32922283644
$ws.Cells.Item(12, 19).Value = 0.03474194025179955
$ws.Cells.Item(12, 20).Value = 0.03756378979244185
$ws.Cells.Item(13, 7).Value = 0.8171586666666667
$ws.Cells.Item(13, 8).Value = 2.451476
$ws.Cells.Item(13, 9).Value = 0.08154160546066344
$ws.Cells.Item(13, 10).Value = 0.08686056894424524
$ws.Cells.Item(13, 15).Value = 0.1743777127077069
$ws.Cells.Item(13, 16).Value = 0.1769957741547643
$ws.Cells.Item(13, 17).Value = 12.20229805814222
$ws.Cells.Item(13, 18).Value = 109.82068252328
$ws.Cells.Item(13, 19).Value = 0.01421903865074475
$ws.Cells.Item(13, 20).Value = 0.01537395364380996
$ws.Cells.Item(14, 7).Value = 0.8171586666666667
$ws.Cells.Item(14, 8).Value = 2.451476
$ws.Cells.Item(14, 9).Value = 0.08154160546066344
$ws.Cells.Item(14, 10).Value = 0.08686056894424524
$ws.Cells.Item(14, 13).Value = 9.680823666666667
$ws.Cells.Item(14, 14).Value = 29.042471
$ws.Cells.Item(14, 15).Value = 0.1130493445068016
$ws.Cells.Item(14, 16).Value = 0.1147466378470605
$ws.Cells.Item(14, 17).Value = 7.910768959688445
$ws.Cells.Item(14, 18).Value = 71.196920637196
$ws.Cells.Item(14, 19).Value = 0.009218225047360237
$ws.Cells.Item(14, 20).Value = 0.009966958247834942
$ws.Cells.Item(15, 7).Value = 0.8171586666666667
$ws.Cells.Item(15, 8).Value = 2.451476
$ws.Cells.Item(15, 9).Value = 0.08154160546066344
$ws.Cells.Item(15, 10).Value = 0.08686056894424524
$ws.Cells.Item(15, 13).Value = 3.79999
$ws.Cells.Item(15, 14).Value = 7.59998
$ws.Cells.Item(15, 15).Value = 0.04437498227672168
$ws.Cells.Item(15, 16).Value = 0.0300274777826206
$ws.Cells.Item(15, 17).Value = 3.105194761746667
$ws.Cells.Item(15, 18).Value = 18.63116857048
$ws.Cells.Item(15, 19).Value = 0.003618407297132372
$ws.Cells.Item(15, 20).Value = 0.002608203804159109
$ws.Cells.Item(16, 7).Value = 0.8171586666666667
$ws.Cells.Item(16, 8).Value = 2.451476
$ws.Cells.Item(16, 9).Value = 0.08154160546066344
$ws.Cells.Item(16, 10).Value = 0.08686056894424524
$ws.Cells.Item(16, 13).Value = 20.734808
$ws.Cells.Item(16, 14).Value = 62.204424
$ws.Cells.Item(16, 15).Value = 0.242133989171347
$ws.Cells.Item(16, 16).Value = 0.245769325661477
$ws.Cells.Item(16, 17).Value = 16.94362805886934
$ws.Cells.Item(16, 18).Value = 152.492652529824
$ws.Cells.Item(16, 19).Value = 0.01974399421362653
$ws.Cells.Item(16, 20).Value = 0.02134766345599938
$ws.Cells.Item(17, 7).Value = 1.8409955
$ws.Cells.Item(17, 8).Value = 3.681991
$ws.Cells.Item(17, 9).Value = 0.1837069529326701
$ws.Cells.Item(17, 10).Value = 0.1304601118296041
$ws.Cells.Item(17, 13).Value = 36.48539666666667
$ws.Cells.Item(17, 14).Value = 109.45619
$ws.Cells.Item(17, 15).Value = 0.4260639713374229
$ws.Cells.Item(17, 16).Value = 0.4324607845540777
$ws.Cells.Item(17, 17).Value = 67.16945107904833
$ws.Cells.Item(17, 18).Value = 403.01670647429
$ws.Cells.Item(17, 19).Value = 0.07827091392879044
$ws.Cells.Item(17, 20).Value = 0.05641888231484329
$ws.Cells.Item(18, 7).Value = 1.8409955
$ws.Cells.Item(18, 8).Value = 3.681991
$ws.Cells.Item(18, 9).Value = 0.1837069529326701
$ws.Cells.Item(18, 10).Value = 0.1304601118296041
$ws.Cells.Item(18, 15).Value = 0.1743777127077069
$ws.Cells.Item(18, 16).Value = 0.1769957741547643
$ws.Cells.Item(18, 17).Value = 27.49083712999667
$ws.Cells.Item(18, 18).Value = 164.94502277998
$ws.Cells.Item(18, 19).Value = 0.03203439826090138
$ws.Cells.Item(18, 20).Value = 0.0230908884895979
$ws.Cells.Item(19, 7).Value = 1.8409955
$ws.Cells.Item(19, 8).Value = 3.681991
$ws.Cells.Item(19, 9).Value = 0.1837069529326701
$ws.Cells.Item(19, 10).Value = 0.1304601118296041
$ws.Cells.Item(19, 13).Value = 9.680823666666667
$ws.Cells.Item(19, 14).Value = 29.042471
$ws.Cells.Item(19, 15).Value = 0.1130493445068016
$ws.Cells.Item(19, 16).Value = 0.1147466378470605
$ws.Cells.Item(19, 17).Value = 17.82235280662684
$ws.Cells.Item(19, 18).Value = 106.934116839761
$ws.Cells.Item(19, 19).Value = 0.02076795061038022
$ws.Cells.Item(19, 20).Value = 0.0149698592055986
$ws.Cells.Item(20, 7).Value = 1.8409955
$ws.Cells.Item(20, 8).Value = 3.681991
$ws.Cells.Item(20, 9).Value = 0.1837069529326701
$ws.Cells.Item(20, 10).Value = 0.1304601118296041
$ws.Cells.Item(20, 13).Value = 3.79999
$ws.Cells.Item(20, 14).Value = 7.59998
$ws.Cells.Item(20, 15).Value = 0.04437498227672168
$ws.Cells.Item(20, 16).Value = 0.0300274777826206
$ws.Cells.Item(20, 17).Value = 6.995764490045
$ws.Cells.Item(20, 18).Value = 27.98305796018
$ws.Cells.Item(20, 19).Value = 0.00815199278049778
$ws.Cells.Item(20, 20).Value = 0.003917388109481636
$ws.Cells.Item(21, 7).Value = 1.8409955
$ws.Cells.Item(21, 8).Value = 3.681991
$ws.Cells.Item(21, 9).Value = 0.1837069529326701
$ws.Cells.Item(21, 10).Value = 0.1304601118296041
$ws.Cells.Item(21, 13).Value = 20.734808
$ws.Cells.Item(21, 14).Value = 62.204424
$ws.Cells.Item(21, 15).Value = 0.242133989171347
$ws.Cells.Item(21, 16).Value = 0.245769325661477
$ws.Cells.Item(21, 17).Value = 38.17268822136401
$ws.Cells.Item(21, 18).Value = 229.036129328184
$ws.Cells.Item(21, 19).Value = 0.0444816973521003
$ws.Cells.Item(21, 20).Value = 0.03206309371008267
$ws.Cells.Item(22, 7).Value = 1.637187333333333
$ws.Cells.Item(22, 8).Value = 4.911562
$ws.Cells.Item(22, 9).Value = 0.1633695988863799
$ws.Cells.Item(22, 10).Value = 0.1740262069565172
$ws.Cells.Item(22, 13).Value = 36.48539666666667
$ws.Cells.Item(22, 14).Value = 109.45619
$ws.Cells.Item(22, 15).Value = 0.4260639713374229
$ws.Cells.Item(22, 16).Value = 0.4324607845540777
$ws.Cells.Item(22, 17).Value = 59.73342927430889
$ws.Cells.Item(22, 18).Value = 537.6008634687801
$ws.Cells.Item(22, 19).Value = 0.06960590009733283
$ws.Cells.Item(22, 20).Value = 0.07525950999338572
$ws.Cells.Item(23, 7).Value = 1.637187333333333
$ws.Cells.Item(23, 8).Value = 4.911562
$ws.Cells.Item(23, 9).Value = 0.1633695988863799
$ws.Cells.Item(23, 10).Value = 0.1740262069565172
$ws.Cells.Item(23, 15).Value = 0.1743777127077069
$ws.Cells.Item(23, 16).Value = 0.1769957741547643
$ws.Cells.Item(23, 17).Value = 24.44745265915111
$ws.Cells.Item(23, 18).Value = 220.02707393236
$ws.Cells.Item(23, 19).Value = 0.02848801697978246
$ws.Cells.Item(23, 20).Value = 0.030801903223486
$ws.Cells.Item(24, 7).Value = 1.637187333333333
$ws.Cells.Item(24, 8).Value = 4.911562
$ws.Cells.Item(24, 9).Value = 0.1633695988863799
$ws.Cells.Item(24, 10).Value = 0.1740262069565172
$ws.Cells.Item(24, 13).Value = 9.680823666666667
$ws.Cells.Item(24, 14).Value = 29.042471
$ws.Cells.Item(24, 15).Value = 0.1130493445068016
$ws.Cells.Item(24, 16).Value = 0.1147466378470605
$ws.Cells.Item(24, 17).Value = 15.84932188330022
$ws.Cells.Item(24, 18).Value = 142.643896949702
$ws.Cells.Item(24, 19).Value = 0.01846882606644435
$ws.Cells.Item(24, 20).Value = 0.01996892214553709
$ws.Cells.Item(25, 7).Value = 1.637187333333333
$ws.Cells.Item(25, 8).Value = 4.911562
$ws.Cells.Item(25, 9).Value = 0.1633695988863799
$ws.Cells.Item(25, 10).Value = 0.1740262069565172
$ws.Cells.Item(25, 13).Value = 3.79999
$ws.Cells.Item(25, 14).Value = 7.59998
$ws.Cells.Item(25, 15).Value = 0.04437498227672168
$ws.Cells.Item(25, 16).Value = 0.0300274777826206
$ws.Cells.Item(25, 17).Value = 6.221295494793334
$ws.Cells.Item(25, 18).Value = 37.32777296876
$ws.Cells.Item(25, 19).Value = 0.007249523055138238
$ws.Cells.Item(25, 20).Value = 0.005225568062980556
$ws.Cells.Item(26, 7).Value = 1.637187333333333
$ws.Cells.Item(26, 8).Value = 4.911562
$ws.Cells.Item(26, 9).Value = 0.1633695988863799
$ws.Cells.Item(26, 10).Value = 0.1740262069565172
$ws.Cells.Item(26, 13).Value = 20.734808
$ws.Cells.Item(26, 14).Value = 62.204424
$ws.Cells.Item(26, 15).Value = 0.242133989171347
$ws.Cells.Item(26, 16).Value = 0.245769325661477
$ws.Cells.Item(26, 17).Value = 33.94676501669867
$ws.Cells.Item(26, 18).Value = 305.520885150288
$ws.Cells.Item(26, 19).Value = 0.03955733268768202
$ws.Cells.Item(26, 20).Value = 0.04277030353112787
